$wb = $excel.ActiveWorkbook

# "2005" sheet: the tube name in A6 was a verbatim duplicate of the other
# sheets' "@terminal" entry; normalize-compare tests need a differently
# cased/accented variant here instead ("@térmiNal") so the migration code's
# name-normalization gets exercised.
$ws2005 = $wb.Worksheets.Item("2005")
$ws2005.Range("A6").Value = "@térmiNal"

# "2011" sheet: move the active selection off this sheet...
$ws2011 = $wb.Worksheets.Item("2011")
$ws2011.Range("A3").Select()

# ...and onto "2005", which becomes the workbook's active/selected tab.
$ws2005.Activate()
$ws2005.Range("E9").Select()
